$d = $word.ActiveDocument
$d.Content.Find.Execute("creating the mailbox, ti returns an error.", $true, $false, $false, $false, $false, $true, 1, $false, "creating the mailbox, it returns an error.", 2)
